# Add a new "as of" column (AF) and a new forecast-target row (44) to both
# the "cases" and "deaths" sheets, matching a new daily data refresh
# (as-of date 2020-05-17).

$wb = $excel.ActiveWorkbook

$sheetNames = @("cases", "deaths")

# Per-sheet data:
#  - b28 / b29: updated "Observed" column values (col B) for rows 28 & 29
#  - b30      : newly-populated "Observed" value for row 30 (was empty)
#  - af       : new diagonal values for column AF (col 32), rows 31..44
#  - af44     : value for the new row 44 / new column AF intersection
$sheetData = @{
    "cases" = @{
        B28 = 91299
        B29 = 96396
        B30 = 101147
        AF = @{
            31 = 109699
            32 = 117719
            33 = 125946
            34 = 133821
            35 = 141186
            36 = 147572
            37 = 153791
            38 = 159822
            39 = 165510
            40 = 170606
            41 = 175511
            42 = 180446
            43 = 184565
            44 = 188615
        }
    }
    "deaths" = @{
        B28 = $null
        B29 = 6724
        B30 = 7025
        AF = @{
            31 = 7590
            32 = 8095
            33 = 8604
            34 = 9086
            35 = 9529
            36 = 9904
            37 = 10277
            38 = 10637
            39 = 10969
            40 = 11267
            41 = 11547
            42 = 11840
            43 = 12070
            44 = 12302
        }
    }
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $data = $sheetData[$name]

    # --- New header cell AF1: same "date" shared string as AE1 would be if
    #     the sequence continued (header strings are column-index-like and
    #     this one already exists in the shared string table, index 31 /
    #     "2020-05-03"). Force text so it isn't auto-parsed into a date
    #     serial, then drop back to the default style so no explicit cell
    #     style (s="..") is left behind.
    $ws.Cells.Item(1, 32).NumberFormat = "@"
    $ws.Cells.Item(1, 32).Value = "2020-05-03"
    $ws.Cells.Item(1, 32).Style = "Normal"

    # --- Rows 2..30: materialize the now-used empty AF cells (sheet's used
    #     range grows to column AF). Setting Style on an empty cell is
    #     enough to persist an empty <c> entry without adding any value.
    $ws.Range("AF2:AF30").Style = "Normal"

    # --- Updated "Observed" (column B) values for rows 28-30.
    if ($null -ne $data.B28) {
        $ws.Cells.Item(28, 2).Value = $data.B28
    }
    if ($null -ne $data.B29) {
        $ws.Cells.Item(29, 2).Value = $data.B29
    }
    if ($null -ne $data.B30) {
        $ws.Cells.Item(30, 2).Value = $data.B30
    }

    # --- New forecast diagonal values in column AF, rows 31-44.
    foreach ($row in $data.AF.Keys) {
        $ws.Cells.Item($row, 32).Value = $data.AF[$row]
    }

    # --- New row 44: label in column A (new date string, forced to text),
    #     empty placeholder cells B44:AE44, and the AF44 value (set above
    #     via the AF loop).
    $ws.Cells.Item(44, 1).NumberFormat = "@"
    $ws.Cells.Item(44, 1).Value = "2020-05-17"
    $ws.Cells.Item(44, 1).Style = "Normal"

    $ws.Range("B44:AE44").Style = "Normal"
}
